# GoodInfo_v2 - 2022-01-19 未完成
# Update row 15 (2022-01-18): B15/C15/E15 switch from text to numeric values,
# and append a new row 16 for 2022-01-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    # Force literal text so Excel does not auto-detect numbers/dates/percentages
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Row 15 (2022-01-18): B15, C15, E15 become real numbers ---
$ws.Range("B15").Value = 54446
$ws.Range("C15").Value = 252
$ws.Range("E15").Value = 0

# --- Row 16 (new): 2022-01-19, 未完成 ---
Set-TextValue $ws.Range("A16") "2022-01-19"
Set-TextValue $ws.Range("B16") "54446.0"
Set-TextValue $ws.Range("C16") "175.0"
Set-TextValue $ws.Range("D16") "0.32%"
Set-TextValue $ws.Range("E16") "0"
Set-TextValue $ws.Range("F16") ""
Set-TextValue $ws.Range("G16") ""

Write-Output "done"
